# Add data for 2022-06-20:
# - Rename sheet "Through 2022-06-11" -> "Through 2022-06-12"
# - Update header label "2022 (through 06-11)" -> "2022 (through 06-12)"
# - Update June "through" total in column I (row 7) 40 -> 47
# - Update grand Total (row 14) in column I 703 -> 710

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Through 2022-06-12"

# Update the "Total" column header text (shared string)
$ws.Range("I1").Value = "2022 (through 06-12)"

# Update June row's "through" total value
$ws.Range("I7").Value = 47

# Update grand total value
$ws.Range("I14").Value = 710
